$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# The inserted row copies formatting from the row above (header), so reset
# it back to normal/default formatting to match the rest of the data rows.
$ws.Range("A2:R2").Style = "Normal"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with data (same constant columns as the rest,
# new values for D, J, K, L, M, N, P, Q)
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 45043
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 100112041
$ws.Range("G2").Value = "Fruto del paraíso"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 24000
$ws.Range("L2").Value = 24000
$ws.Range("M2").Value = 24000
$ws.Range("N2").Value = "$/caja 18 kilos empedrada"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 1333
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = "Hortaliza"
